$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

# Require angle brackets for @base and @prefix values (column D, rows 1-3)
$ws.Range("D1").Value = "<http://sales.data/purchases/2015>"
$ws.Range("D2").Value = "<http://sales.data/purchases#>"
$ws.Range("D3").Value = "<http://sales.data/schema#>"

# Update the active selection to the edited range
$ws.Range("D1:D3").Select()
